# Adding more noise results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The per-row "Faults" differences (column G) are no longer needed -
# drop the whole helper column.
$ws.Columns.Item(7).Delete()

# Drop the stray trailing blank row left over at the bottom of the sheet.
$ws.Rows.Item(9).Delete()

# New "Partial Occlusion" section.
$ws.Range("A5").Value = "Partial Occlusion"
$ws.Range("C5").Value = 500
$ws.Range("D5").Value = 67
$ws.Range("E5").Value = 433
$ws.Range("F5").Value = 393

# New "Noise" section, broken down by noise type.
$ws.Range("A6").Value = "Noise"
$ws.Range("B6").Value = "Gaussian"
$ws.Range("C6").Value = 24
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 23
$ws.Range("F6").Value = 21

$ws.Range("B7").Value = "S&P"
$ws.Range("C7").Value = 24
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 22
$ws.Range("F7").Value = 22

$ws.Range("B8").Value = "Poisson"
$ws.Range("C8").Value = 24
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 24
$ws.Range("F8").Value = 23

# Update the active selection / view state to reflect where editing left off.
$ws.Range("C9").Select() | Out-Null

$w = $wb.Windows.Item(1)
$w.TabRatio = 0.993
